$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to machine-friendly snake_case identifiers
$ws.Range("A1").Value = "tipo_identificacion"
$ws.Range("B1").Value = "numero_identificacion"
$ws.Range("C1").Value = "primer_apellido"
$ws.Range("D1").Value = "segundo_apellido"
$ws.Range("E1").Value = "primer_nombre"
$ws.Range("F1").Value = "segundo_nombre"
$ws.Range("G1").Value = "regional"
$ws.Range("H1").Value = "fecha_gestion"
$ws.Range("I1").Value = "nombre"
$ws.Range("J1").Value = "ciex"
$ws.Range("K1").Value = "medico_id"

# Shorter labels no longer need as much vertical space
$ws.Rows.Item(1).RowHeight = 26.4

# Update the active selection shown when the workbook is opened
$ws.Range("K5").Select()
